$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function RoundSig16($x) {
    # Matches the source data pipeline's float formatting: round to 16
    # significant decimal digits (as opposed to 16 decimal places).
    return [double]($x.ToString("G16"))
}

# Row 16 (cohort 2019, period 6): num_customers 9 -> 10, retention_rate recalculated
$ws.Range("C16").Value = 10
$ws.Range("E16").Value = RoundSig16 ($ws.Range("C16").Value2 / $ws.Range("D16").Value2)

# Row 27 (cohort 2021, period 4): num_customers 48 -> 49, retention_rate recalculated
$ws.Range("C27").Value = 49
$ws.Range("E27").Value = RoundSig16 ($ws.Range("C27").Value2 / $ws.Range("D27").Value2)

# Row 34 (cohort 2023, period 2): num_customers 84 -> 85, retention_rate recalculated
$ws.Range("C34").Value = 85
$ws.Range("E34").Value = RoundSig16 ($ws.Range("C34").Value2 / $ws.Range("D34").Value2)

# Row 36 (cohort 2024, period 1): num_customers 140 -> 142, retention_rate recalculated
$ws.Range("C36").Value = 142
$ws.Range("E36").Value = RoundSig16 ($ws.Range("C36").Value2 / $ws.Range("D36").Value2)

# Row 37 (cohort 2025, period 0): num_customers 895 -> 904, cohort_size 895 -> 904, retention_rate stays 1
$ws.Range("C37").Value = 904
$ws.Range("D37").Value = 904
$ws.Range("E37").Value = RoundSig16 ($ws.Range("C37").Value2 / $ws.Range("D37").Value2)
